$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.63849413394928
$ws.Range("B1").Value = 1.531911015510559
$ws.Range("C1").Value = 4.803493022918701
$ws.Range("D1").Value = 1.244050621986389
$ws.Range("E1").Value = 0.6421604156494141
